# Regenerate the "K" column (column G) values on the active worksheet.
#
# Context (from commit message): "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals" — the K column
# (strikeouts) was recalculated/rewritten with new values for every
# existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 0
    22 = 1
    23 = 1
    24 = 3
    25 = 2
    26 = 1
    27 = 2
    28 = 3
    29 = 1
    30 = 0
    31 = 2
    32 = 2
    33 = 2
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 0
    41 = 3
    42 = 3
    43 = 2
    44 = 1
    46 = 1
    48 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
